$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 236.7946826666667
$ws.Range("H2").Value = 710.384048
$ws.Range("I2").Value = 0.7123899543147419
$ws.Range("J2").Value = 0.7240508783182559
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 28872.07705945182
$ws.Range("R2").Value = 259848.6935350664
$ws.Range("S2").Value = 0.1625839258149097
$ws.Range("T2").Value = 0.1752000703581473
$ws.Range("G3").Value = 236.7946826666667
$ws.Range("H3").Value = 710.384048
$ws.Range("I3").Value = 0.7123899543147419
$ws.Range("J3").Value = 0.7240508783182559
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 35025.32683420261
$ws.Range("R3").Value = 315227.9415078235
$ws.Range("S3").Value = 0.197233996290916
$ws.Range("T3").Value = 0.2125389078532029
$ws.Range("G4").Value = 236.7946826666667
$ws.Range("H4").Value = 710.384048
$ws.Range("I4").Value = 0.7123899543147419
$ws.Range("J4").Value = 0.7240508783182559
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 19773.5327143764
$ws.Range("R4").Value = 177961.7944293876
$ws.Range("S4").Value = 0.1113483650418707
$ws.Range("T4").Value = 0.1199887460695787
$ws.Range("G5").Value = 236.7946826666667
$ws.Range("H5").Value = 710.384048
$ws.Range("I5").Value = 0.7123899543147419
$ws.Range("J5").Value = 0.7240508783182559
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 21564.52767862874
$ws.Range("R5").Value = 129387.1660717724
$ws.Range("S5").Value = 0.121433783967682
$ws.Range("T5").Value = 0.08723784711334985
$ws.Range("G6").Value = 236.7946826666667
$ws.Range("H6").Value = 710.384048
$ws.Range("I6").Value = 0.7123899543147419
$ws.Range("J6").Value = 0.7240508783182559
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 21272.59949801008
$ws.Range("R6").Value = 191453.3954820907
$ws.Range("S6").Value = 0.1197898831993635
$ws.Range("T6").Value = 0.1290853069239771
$ws.Range("I7").Value = 0.2358656137148928
$ws.Range("J7").Value = 0.2397264359793184
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 9559.2731671825
$ws.Range("R7").Value = 86033.45850464249
$ws.Range("S7").Value = 0.05383000870555188
$ws.Range("T7").Value = 0.05800709550665497
$ws.Range("I8").Value = 0.2358656137148928
$ws.Range("J8").Value = 0.2397264359793184
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 11596.55629515513
$ws.Range("S8").Value = 0.06530232114986338
$ws.Range("T8").Value = 0.07036963342262467
$ws.Range("I9").Value = 0.2358656137148928
$ws.Range("J9").Value = 0.2397264359793184
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 6546.830708013235
$ws.Range("R9").Value = 58921.47637211911
$ws.Range("S9").Value = 0.03686639641348359
$ws.Range("T9").Value = 0.03972714530738287
$ws.Range("I10").Value = 0.2358656137148928
$ws.Range("J10").Value = 0.2397264359793184
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 7139.812296039705
$ws.Range("R10").Value = 42838.87377623823
$ws.Range("S10").Value = 0.04020558376459733
$ws.Range("T10").Value = 0.02888363069121179
$ws.Range("I11").Value = 0.2358656137148928
$ws.Range("J11").Value = 0.2397264359793184
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 7043.157621075169
$ws.Range("R11").Value = 63388.41858967653
$ws.Range("S11").Value = 0.03966130368139654
$ws.Range("T11").Value = 0.04273893105144407
$ws.Range("G12").Value = 0.6305213333333334
$ws.Range("H12").Value = 1.891564
$ws.Range("I12").Value = 0.001896905195629352
$ws.Range("J12").Value = 0.001927955138422806
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 76.878671085932
$ws.Range("R12").Value = 691.9080397733879
$ws.Range("S12").Value = 0.0004329178025829682
$ws.Range("T12").Value = 0.0004665112439108973
$ws.Range("G13").Value = 0.6305213333333334
$ws.Range("H13").Value = 1.891564
$ws.Range("I13").Value = 0.001896905195629352
$ws.Range("J13").Value = 0.001927955138422806
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 93.26314057070665
$ws.Range("R13").Value = 839.3682651363599
$ws.Range("S13").Value = 0.0005251817351619786
$ws.Range("T13").Value = 0.0005659346487668257
$ws.Range("G14").Value = 0.6305213333333334
$ws.Range("H14").Value = 1.891564
$ws.Range("I14").Value = 0.001896905195629352
$ws.Range("J14").Value = 0.001927955138422806
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 52.65166460401245
$ws.Range("R14").Value = 473.864981436112
$ws.Range("S14").Value = 0.0002964911154255833
$ws.Range("T14").Value = 0.0003194981547085029
$ws.Range("G15").Value = 0.6305213333333334
$ws.Range("H15").Value = 1.891564
$ws.Range("I15").Value = 0.001896905195629352
$ws.Range("J15").Value = 0.001927955138422806
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 57.42060840011668
$ws.Range("R15").Value = 344.5236504007
$ws.Range("S15").Value = 0.0003233459067440159
$ws.Range("T15").Value = 0.0002322912113548987
$ws.Range("G16").Value = 0.6305213333333334
$ws.Range("H16").Value = 1.891564
$ws.Range("I16").Value = 0.001896905195629352
$ws.Range("J16").Value = 0.001927955138422806
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 56.64328120844
$ws.Range("R16").Value = 509.78953087596
$ws.Range("S16").Value = 0.0003189686357148053
$ws.Range("T16").Value = 0.0003437198796816814
$ws.Range("G17").Value = 16.059769
$ws.Range("H17").Value = 32.119538
$ws.Range("I17").Value = 0.04831535056182164
$ws.Range("J17").Value = 0.032737474561192
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 1958.147382801291
$ws.Range("R17").Value = 11748.88429680774
$ws.Range("S17").Value = 0.01102668464636154
$ws.Range("T17").Value = 0.007921553606551687
$ws.Range("G18").Value = 16.059769
$ws.Range("H18").Value = 32.119538
$ws.Range("I18").Value = 0.04831535056182164
$ws.Range("J18").Value = 0.032737474561192
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 2375.46997158977
$ws.Range("R18").Value = 14252.81982953862
$ws.Range("S18").Value = 0.01337670417134269
$ws.Range("T18").Value = 0.009609804086239066
$ws.Range("G19").Value = 16.059769
$ws.Range("H19").Value = 32.119538
$ws.Range("I19").Value = 0.04831535056182164
$ws.Range("J19").Value = 0.032737474561192
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 1341.070517845417
$ws.Range("R19").Value = 8046.423107072504
$ws.Range("S19").Value = 0.007551812401199015
$ws.Range("T19").Value = 0.0054252106305098
$ws.Range("G20").Value = 16.059769
$ws.Range("H20").Value = 32.119538
$ws.Range("I20").Value = 0.04831535056182164
$ws.Range("J20").Value = 0.032737474561192
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 1462.538470935163
$ws.Range("R20").Value = 5850.153883740651
$ws.Range("S20").Value = 0.008235820573987087
$ws.Range("T20").Value = 0.003944400712944262
$ws.Range("G21").Value = 16.059769
$ws.Range("H21").Value = 32.119538
$ws.Range("I21").Value = 0.04831535056182164
$ws.Range("J21").Value = 0.032737474561192
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 1442.73946576947
$ws.Range("R21").Value = 8656.43679461682
$ws.Range("S21").Value = 0.008124328768931302
$ws.Range("T21").Value = 0.005836505524947183
$ws.Range("G22").Value = 0.5092873333333333
$ws.Range("H22").Value = 1.527862
$ws.Range("I22").Value = 0.001532176212914103
$ws.Range("J22").Value = 0.001557256002810873
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 62.096762342006
$ws.Range("R22").Value = 558.870861078054
$ws.Range("S22").Value = 0.0003496781814889789
$ws.Range("T22").Value = 0.000376812416679632
$ws.Range("G23").Value = 0.5092873333333333
$ws.Range("H23").Value = 1.527862
$ws.Range("I23").Value = 0.001532176212914103
$ws.Range("J23").Value = 0.001557256002810873
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 75.33089468748665
$ws.Range("R23").Value = 677.9780521873799
$ws.Range("S23").Value = 0.0004242019917105902
$ws.Range("T23").Value = 0.0004571191058479543
$ws.Range("G24").Value = 0.5092873333333333
$ws.Range("H24").Value = 1.527862
$ws.Range("I24").Value = 0.001532176212914103
$ws.Range("J24").Value = 0.001557256002810873
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 42.52802315185511
$ws.Range("R24").Value = 382.752208366696
$ws.Range("S24").Value = 0.0002394830460911513
$ws.Range("T24").Value = 0.0002580663882634913
$ws.Range("G25").Value = 0.5092873333333333
$ws.Range("H25").Value = 1.527862
$ws.Range("I25").Value = 0.001532176212914103
$ws.Range("J25").Value = 0.001557256002810873
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 46.38001441739167
$ws.Range("R25").Value = 278.28008650435
$ws.Range("S25").Value = 0.0002611743106602396
$ws.Range("T25").Value = 0.0001876272305685233
$ws.Range("G26").Value = 0.5092873333333333
$ws.Range("H26").Value = 1.527862
$ws.Range("I26").Value = 0.001532176212914103
$ws.Range("J26").Value = 0.001557256002810873
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 45.75214844102
$ws.Range("R26").Value = 411.76933596918
$ws.Range("S26").Value = 0.0002576386829631426
$ws.Range("T26").Value = 0.0002776308614512716
